$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(57, 1).Value = 'li: Cheltenham East Primary School, CHELTENHAM'
$ws.Cells.Item(67, 1).Value = 'li: PCW, WINDSOR'
$ws.Cells.Item(68, 1).Value = 'li: Thomas Mitchell Primary School, ENDEAVOUR HILLS'
$ws.Cells.Item(69, 1).Value = 'li: Westall Primary School, CLAYTON SOUTH'
$ws.Cells.Item(70, 1).Value = 'TAFEThe Department hasnotbeen advised of any TAFE closures.North-Eastern Victoria RegionEarly childhood services'
$ws.Cells.Item(71, 1).Value = 'li: The Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(72, 1).Value = 'li: Aqualink Box Hill Creche BOX HILL'
$ws.Cells.Item(73, 1).Value = 'li: Aqualink Nunawading Creche FOREST HILL'
$ws.Cells.Item(74, 1).Value = 'li: Camp Australia - Kerrimuir Primary School OSHC BOX HILL NORTH'
$ws.Cells.Item(75, 1).Value = 'li: Camp Australia - Mount Waverley Primary School OSHC MOUNT WAVERLEY'
$ws.Cells.Item(76, 1).Value = 'li: Camp Australia - Mountain Gate Primary School OSHC FERNTREE GULLY'
$ws.Cells.Item(77, 1).Value = 'li: Camp Australia - Scoresby Primary School OSHC SCORESBY'
$ws.Cells.Item(78, 1).Value = 'li: Camp Australia - St Clement of Rome School OSHC BULLEEN'
$ws.Cells.Item(79, 1).Value = 'li: Camp Australia - St Timothy''s School Vermont OSHC VERMONT'
$ws.Cells.Item(80, 1).Value = 'li: Camp Australia - Templestowe Heights Primary School OSHC TEMPLESTOWE LOWER'
$ws.Cells.Item(81, 1).Value = 'li: Carey Donvale OSHClub DONVALE'
$ws.Cells.Item(82, 1).Value = 'li: Clever Kids Childcare - Ashburton ASHBURTON'
$ws.Cells.Item(83, 1).Value = 'li: Fitness First Doncaster (Playzone) DONCASTER'
$ws.Cells.Item(84, 1).Value = 'li: Flamingo Community Group WANTIRNA SOUTH'
$ws.Cells.Item(85, 1).Value = 'li: Hawthorn Early Years HAWTHORN'
$ws.Cells.Item(86, 1).Value = 'li: MakerDojo HAWTHORN'
$ws.Cells.Item(87, 1).Value = 'li: Paisley Park Early Learning Centre Chadstone HUGHESDALE'
$ws.Cells.Item(88, 1).Value = 'li: St Andrews Christian College Outside School Hours Care WANTIRNA SOUTH'
$ws.Cells.Item(89, 1).Value = 'li: Vermont Primary School Kindergarten VERMONT'
$ws.Cells.Item(90, 1).Value = 'li: Vermont Primary School Outside School Hours Child Care Service VERMONT'
$ws.Cells.Item(91, 1).Value = 'Schools closedThe Department hasbeen advisedof the followingschool closures:Aquinas College, RINGWOODAshwood High School, BURWOODAuburn High School, HAWTHORN EASTGlen Iris Primary School, GLEN IRISKerrimuir Primary School, BOX HILL NORTH'
$ws.Cells.Item(92, 1).Value = 'li: Mooroolbark College, MOOROOLBARK'
$ws.Cells.Item(93, 1).Value = 'Mount Waverley Primary School, MOUNT WAVERLEYOur Lady Of the Way, WALLAN EASTVermont Primary School, VERMONTWattle Park Primary School, BURWOODTAFE'
$ws.Cells.Item(94, 1).Value = 'li: The Department has'
$ws.Cells.Item(95, 1).Value = 'not'
$ws.Cells.Item(96, 1).Value = 'li: been advised of any TAFE closures.'
$ws.Cells.Item(97, 1).Value = 'North-Western Victoria RegionEarly childhood services'
$ws.Cells.Item(98, 1).Value = 'li: The Department has been advised of the following early childhood service closures:'
$ws.Cells.Item(99, 1).Value = 'li: Big Childcare - Fitzroy PS OSHC FITZROY'
$ws.Cells.Item(100, 1).Value = 'li: Bright Stars Early Years Child Care Centre EPPING'
$ws.Cells.Item(101, 1).Value = 'li: Camp Australia - Ivanhoe Grammar OSHC IVANHOE'
$ws.Cells.Item(102, 1).Value = 'li: Camp Australia - Ivanhoe Grammar Plenty Valley Campus OSHC MERNDA'
$ws.Cells.Item(103, 1).Value = 'li: Camp Australia - St Joseph''s Primary School - Mernda OSHC MERNDA'
$ws.Cells.Item(104, 1).Value = 'li: Camp Australia - Thomastown Meadows Primary School OSHC THOMASTOWN'
$ws.Cells.Item(105, 1).Value = 'li: Collingwood College Afterschool Care and Vacation Care Program COLLINGWOOD'
$ws.Cells.Item(106, 1).Value = 'li: Diamond Creek Community Centre DIAMOND CREEK'
$ws.Cells.Item(107, 1).Value = 'li: Diamond Valley Sports and Fitness Centre GREENSBOROUGH'
$ws.Cells.Item(108, 1).Value = 'li: Eltham Leisure Centre ELTHAM'
$ws.Cells.Item(109, 1).Value = 'li: Greenvale PS TheirCare GREENVALE'
$ws.Cells.Item(110, 1).Value = 'li: Kangaroo Ground Primary Combined OSHC KANGAROO GROUND'
$ws.Cells.Item(111, 1).Value = 'li: New Futures Broadmeadows BROADMEADOWS'
$ws.Cells.Item(112, 1).Value = 'li: New Futures Epping EPPING'
$ws.Cells.Item(113, 1).Value = 'li: Nino Early Learning Adventures - Bundoora BUNDOORA'
$ws.Cells.Item(114, 1).Value = 'li: Pender''s Grove Primary School Combined OSHC THORNBURY'
$ws.Cells.Item(115, 1).Value = 'li: Richmond West Afterschool Care and Vacation Care Program RICHMOND'
$ws.Cells.Item(116, 1).Value = 'li: St Bernard''s Out of School Hours Care COBURG EAST'
$ws.Cells.Item(117, 1).Value = 'li: Tullamarine Early Learning Centre TULLAMARINE'
$ws.Cells.Item(118, 1).Value = 'Schoolsclosed'
$ws.Cells.Item(119, 1).Value = 'li: The Department has'
$ws.Cells.Item(120, 1).Value = 'li: been advised of the following school closures:'
$ws.Cells.Item(121, 1).Value = 'Al Siraat College, EPPINGCharles La Trobe, MACLEOD WESTEpping Secondary College, EPPINGEpping Views Primary School, EPPINGFitzroy Primary School, FITZROYGisborne Secondary College, GISBORNEGladstone Park Secondary College, GLADSTONE PARKGreenvalePrimary School, GREENVALELalor Secondary College, LALORMelbourne Girls College, RICHMONDPascoe Vale Girls Secondary College, OAK PARKPenders Grove Primary School, THORNBURYPeter Lalor Secondary College, LALORPrinces Hill Secondary College, PRINCES HILLRoxburgh College, ROXBURGH PARKSacred Heart School, FITZROYSt Joseph''s Primary School, MERNDAThornbury High School, THORNBURYTAFE'
$ws.Cells.Item(122, 1).Value = 'li: The Department has'
$ws.Cells.Item(123, 1).Value = 'not'
$ws.Cells.Item(124, 1).Value = 'li: been advised of any TAFE closures.'
$ws.Cells.Item(125, 1).Value = 'South-Western Victoria RegionThe Department has not been advised of any school, early childhood service or TAFE closures, or buscancellations.Early childhood services'
$ws.Cells.Item(126, 1).Value = 'li: The Department has been advised of the followi'
$ws.Cells.Item(127, 1).Value = 'ng early childhood service closures:'
$ws.Cells.Item(128, 1).Value = 'li: Aerotots Activity Centre WERRIBEE'
$ws.Cells.Item(129, 1).Value = 'li: Altona Meadows Community Centre Occasional Care ALTONA MEADOWS'
$ws.Cells.Item(130, 1).Value = 'li: Aquapulse Creche HOPPERS CROSSING'
$ws.Cells.Item(131, 1).Value = 'li: Big Childcare - Keilor PS OSHC KEILOR'
$ws.Cells.Item(132, 1).Value = 'li: Big Childcare - Manor Lakes P-12 College OSHC WYNDHAM VALE'
$ws.Cells.Item(133, 1).Value = 'li: Big Childcare - Overnewton Anglican Community College OSHC TAYLORS LAKES'
$ws.Cells.Item(134, 1).Value = 'li: Big Childcare - Sunshine Heights PS OSHC SUNSHINE'
$ws.Cells.Item(135, 1).Value = 'li: Big Childcare – Sydenham/Hillside Sydenham Campus SYDENHAM'
$ws.Cells.Item(136, 1).Value = 'li: Blackwood Street Neighbourhood House YARRAVILLE'
$ws.Cells.Item(137, 1).Value = 'li: Bluewater Leisure Centre Creche COLAC'
$ws.Cells.Item(138, 1).Value = 'li: Camp Australia - Baden Powell P-9 College Derrimut Heath Campus OSHC HOPPERS CROSSING'
$ws.Cells.Item(139, 1).Value = 'li: Camp Australia - Footscray City Primary School OSHC FOOTSCRAY'
$ws.Cells.Item(140, 1).Value = 'li: Camp Australia - Haileybury City Campus OSHC WEST MELBOURNE'
$ws.Cells.Item(141, 1).Value = 'li: Camp Australia - Kardinia International College OSHC BELL POST HILL'
$ws.Cells.Item(142, 1).Value = 'li: Camp Australia - Melton Christian College OSHC Melton South'
$ws.Cells.Item(143, 1).Value = 'li: Cana Catholic Primary OSHClub HILLSIDE'
$ws.Cells.Item(144, 1).Value = 'li: Carranballac Jamieson OSHClub POINT COOK'
$ws.Cells.Item(145, 1).Value = 'li: Coragulac & District Kindergarten CORAGULAC'
$ws.Cells.Item(146, 1).Value = 'li: Eagle Stadium WERRIBEE'
$ws.Cells.Item(147, 1).Value = 'li: Early Learning Centre Rose Grange TARNEIT'
$ws.Cells.Item(148, 1).Value = 'li: Energy Force Fitness Creche DRYSDALE'
$ws.Cells.Item(149, 1).Value = 'li: Fernwood Fitness Sydenham SYDENHAM'
$ws.Cells.Item(150, 1).Value = 'li: Fun 4 All Occasional Care Center WERRIBEE'
$ws.Cells.Item(151, 1).Value = 'li: Future Kids Child Care-West Tarneit TARNEIT'
$ws.Cells.Item(152, 1).Value = 'li: Genesis Maidstone MAIDSTONE'
$ws.Cells.Item(153, 1).Value = 'li: Goodlife Essendon Child Minding ESSENDON'
$ws.Cells.Item(154, 1).Value = 'li: Goodlife Geelong BELMONT'
$ws.Cells.Item(155, 1).Value = 'li: Goodlife Point Cook POINT COOK'
$ws.Cells.Item(156, 1).Value = 'li: Goodlife Taylors Lakes TAYLORS LAKES'
$ws.Cells.Item(157, 1).Value = 'li: Happy Feet ELC MELTON WEST'
$ws.Cells.Item(158, 1).Value = 'li: Highpoint Kinder Haven MARIBYRNONG'
$ws.Cells.Item(159, 1).Value = 'li: Kardinia International College Kindergarten BELL POST HILL'
$ws.Cells.Item(160, 1).Value = 'li: Keilor Basketball Stadium Creche KEILOR PARK'
$ws.Cells.Item(161, 1).Value = 'li: Kensington Neighbourhood House Inc KENSINGTON'
$ws.Cells.Item(162, 1).Value = 'li: Kids on Collins MELBOURNE'
$ws.Cells.Item(163, 1).Value = 'li: Maribyrnong Aquatic Centre Occasional Child Care MARIBYRNONG'
$ws.Cells.Item(164, 1).Value = 'li: New Futures Braybrook BRAYBROOK'
$ws.Cells.Item(165, 1).Value = 'li: Shuter Street Occasional Care MOONEE PONDS'
$ws.Cells.Item(166, 1).Value = 'li: South Kingsville Community Centre SOUTH KINGSVILLE'
$ws.Cells.Item(167, 1).Value = 'li: Story House Early Learning Keilor Downs KEILOR DOWNS'
$ws.Cells.Item(168, 1).Value = 'li: Sunshine Leisure Centre SUNSHINE'
$ws.Cells.Item(169, 1).Value = 'li: Willaura Primary School OSHC ARARAT'
$ws.Cells.Item(170, 1).Value = 'li: Woodlea Early Education AINTREE'
$ws.Cells.Item(171, 1).Value = 'li: Yarraville Community Centre YARRAVILLE'
$ws.Cells.Item(172, 1).Value = 'Schools closed'
$ws.Cells.Item(173, 1).Value = 'li: The Department has been advised of the following school closures:'
$ws.Cells.Item(174, 1).Value = 'Al Taqwa College, TRUGANINABelmont High School, BELMONTBraybrook College, BRAYBROOKBuckley Park College, ESSENDONCana Primary School, HILLSIDECaroline Chilsholm Catholic College, BRAYBROOKCatholic Regional College, SYDENHAMClonard Secondary College, GEELONG WESTColac Secondary College, COLACCopperfield College Delahey Campus, DELAHEYCopperfield College, Sydenham Campus, SYDENHAM'
$ws.Cells.Item(175, 1).Value = 'li: Footscray High School, FOOTSCRAY'
$ws.Cells.Item(176, 1).Value = 'li: Geelong Baptist College, LOVELY BANKS'
$ws.Cells.Item(177, 1).Value = 'Grovedale West Primary School, GROVEDALEHoppers Crossing Secondary College, HOPPERS CROSSINGLowther Hall Anglican Grammar School, ESSENDONKeilor Downs Secondary College, KEILOR DOWNSKeilor Primary School, KEILORMarian College, SUNSHINE WESTManor Lakes P-12 College, WYNDHAM VALEMacKillop College, WERRIBEEMelton Secondary College, MELTONMount St Joseph''s Girls'' College, ALTONAPoint Cook Senior Secondary College, POINT COOKRosamond Specialist School, BRAYBROOKSt Alban''s Secondary College, ST ALBANSSt Joseph''s Flexible Learning Centre Geelong Campus, GEELONGStaughton College, MELTON SOUTHSunshine Heights Primary School, SUNSHINESydenham Hillside Primary School, SYDENHAMTarneit Rise Primary School, TARNEITTarneit Senior College, TARNEITTaylors Lakes Secondary College, TAYLORS LAKESThomas Carr College, TARNEITVictorian College of the Arts, ALBERT PARKWarracknabeal Secondary College, WARRACKNABEALTAFE'
